$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 252.93333
$ws.Range("I2").Value = 219.16667
$ws.Range("K2").Value = 219.16667
$ws.Range("M2").Value = -106.16667
# Row 9
$ws.Range("H9").Value = 621968.0600000001
$ws.Range("I9").Value = 932886.7
$ws.Range("K9").Value = 932886.7
$ws.Range("M9").Value = -932717.7
# Row 40
$ws.Range("H40").Value = 8155.3335
$ws.Range("I40").Value = 8732
$ws.Range("K40").Value = 8732
$ws.Range("M40").Value = -8557
# Row 64
$ws.Range("H64").Value = 3500
$ws.Range("I64").Value = 3500
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 3500
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -3252
$ws.Range("N64").Value = -3996
# Row 67
$ws.Range("H67").Value = 3500
$ws.Range("I67").Value = 3500
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 3500
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -2642
$ws.Range("N67").Value = -5216
# Row 88
$ws.Range("I88").Value = 12501011
$ws.Range("J88").Value = 1692.6666
$ws.Range("K88").Value = 12501011
$ws.Range("L88").Value = 1692.6666
$ws.Range("M88").Value = -12500605
$ws.Range("N88").Value = -2504.6666
# Row 91
$ws.Range("I91").Value = 12501011
$ws.Range("J91").Value = 1692.6666
$ws.Range("K91").Value = 12501011
$ws.Range("L91").Value = 1692.6666
$ws.Range("M91").Value = -12499607
$ws.Range("N91").Value = -4500.6666
# Row 123
$ws.Range("H123").Value = 109985.734
$ws.Range("J123").Value = 109985.734
$ws.Range("L123").Value = 109985.734
$ws.Range("N123").Value = -119785.734
# Row 129
$ws.Range("H129").Value = 2142.353
$ws.Range("I129").Value = 1402.091
$ws.Range("K129").Value = 4206.272999999999
$ws.Range("M129").Value = 793.7270000000008
# Row 132
$ws.Range("H132").Value = 5915.1113
$ws.Range("I132").Value = 6214.091
$ws.Range("K132").Value = 18642.273
$ws.Range("M132").Value = -16112.273
# Row 138
$ws.Range("H138").Value = 6073.548
$ws.Range("J138").Value = 6816.697
$ws.Range("L138").Value = 20450.091
$ws.Range("N138").Value = -30730.091
# Row 141
$ws.Range("H141").Value = 2741.5
$ws.Range("I141").Value = 2741.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8224.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3044.5
$ws.Range("N141").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 82
$ws.Range("H82").Value = 50181
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
# Row 85
$ws.Range("H85").Value = 50181
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
# Row 139
$ws.Range("H139").Value = 69985.92
$ws.Range("J139").Value = 69985.92
$ws.Range("L139").Value = 69985.92
$ws.Range("N139").Value = -80265.92

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 21
$ws.Range("H21").Value = 11500
$ws.Range("J21").Value = 11500
$ws.Range("L21").Value = 11500
$ws.Range("N21").Value = -11972
# Row 114
$ws.Range("H114").Value = 29999
$ws.Range("J114").Value = 29999
$ws.Range("L114").Value = 29999
$ws.Range("N114").Value = -38677
# Row 115
$ws.Range("H115").Value = 29684
$ws.Range("J115").Value = 29684
$ws.Range("L115").Value = 29684
$ws.Range("N115").Value = -32818
# Row 134
$ws.Range("H134").Value = 2972.7693
$ws.Range("I134").Value = 2487.5
$ws.Range("J134").Value = 3188.4443
$ws.Range("K134").Value = 7462.5
$ws.Range("L134").Value = 9565.332900000001
$ws.Range("M134").Value = -4927.5
$ws.Range("N134").Value = -14635.3329

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 479.33334
$ws.Range("I22").Value = 312.6
$ws.Range("J22").Value = 687.75
$ws.Range("K22").Value = 312.6
$ws.Range("L22").Value = 687.75
$ws.Range("M22").Value = 37.39999999999998
$ws.Range("N22").Value = -1387.75
# Row 58
$ws.Range("H58").Value = 5984.4
$ws.Range("I58").Value = 3810.875
$ws.Range("K58").Value = 3810.875
$ws.Range("M58").Value = -3607.875
# Row 132
$ws.Range("H132").Value = 3346.3057
$ws.Range("I132").Value = 3203.4614
$ws.Range("K132").Value = 9610.3842
$ws.Range("M132").Value = -7080.3842
# Row 136
$ws.Range("H136").Value = 5984.4
$ws.Range("I136").Value = 3810.875
$ws.Range("K136").Value = 11432.625
$ws.Range("M136").Value = -8882.625

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 3800
$ws.Range("I63").Value = 1900
$ws.Range("K63").Value = 5700
$ws.Range("M63").Value = -4951
# Row 66
$ws.Range("H66").Value = 3800
$ws.Range("I66").Value = 1900
$ws.Range("K66").Value = 17100
$ws.Range("M66").Value = -13356
# Row 122
$ws.Range("H122").Value = 1835.6666
$ws.Range("I122").Value = 1321.875
$ws.Range("K122").Value = 11896.875
$ws.Range("M122").Value = -9446.875
# Row 132
$ws.Range("H132").Value = 5178.216
$ws.Range("J132").Value = 5106.049
$ws.Range("L132").Value = 45954.441
$ws.Range("N132").Value = -51014.441
# Row 137
$ws.Range("H137").Value = 7478.4443
$ws.Range("I137").Value = 9884.5
$ws.Range("J137").Value = 2666.3333
$ws.Range("K137").Value = 29653.5
$ws.Range("L137").Value = 7998.999899999999
$ws.Range("M137").Value = -24553.5
$ws.Range("N137").Value = -18198.9999

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 58827300
$ws.Range("J80").Value = 5799.75
$ws.Range("L80").Value = 5799.75
$ws.Range("N80").Value = -7795.75
# Row 83
$ws.Range("H83").Value = 58827300
$ws.Range("J83").Value = 5799.75
$ws.Range("L83").Value = 28998.75
$ws.Range("N83").Value = -38982.75
# Row 132
$ws.Range("H132").Value = 2559.8438
$ws.Range("J132").Value = 2432.9167
$ws.Range("L132").Value = 7298.750100000001
$ws.Range("N132").Value = -12358.7501

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 18
$ws.Range("H18").Value = 15000
$ws.Range("I18").Value = 15000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -14828
$ws.Range("N18").ClearContents()
# Row 40
$ws.Range("H40").Value = 83393.47
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 132
$ws.Range("H132").Value = 7303.647
$ws.Range("I132").Value = 8264.111000000001
$ws.Range("J132").Value = 6223.125
$ws.Range("K132").Value = 24792.333
$ws.Range("L132").Value = 18669.375
$ws.Range("M132").Value = -22262.333
$ws.Range("N132").Value = -23729.375
# Row 140
$ws.Range("H140").Value = 98254
$ws.Range("J140").Value = 98254
$ws.Range("L140").Value = 98254
$ws.Range("N140").Value = -108614

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 48
$ws.Range("H48").Value = 8000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 8000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 8000
$ws.Range("N48").Value = -9138
$ws.Range("M48").ClearContents()
# Row 132
$ws.Range("H132").Value = 2461.2693
$ws.Range("I132").Value = 2370.149
$ws.Range("J132").Value = 3317.8
$ws.Range("K132").Value = 7110.447
$ws.Range("L132").Value = 9953.400000000001
$ws.Range("M132").Value = -4580.447
$ws.Range("N132").Value = -15013.4
